# added 4wk low sales check
# Updates forecast figures (MyForecast, Inventory Coverage, Stockout Risk,
# Reorder Urgency, Seasonality Index) on the "Forecast Comparison" sheet and
# the derived totals on the "Summary" sheet.

$wb = $excel.ActiveWorkbook

$fc = $wb.Worksheets.Item("Forecast Comparison")

# Row 2 (W10)
$fc.Range("D2").Value = 59
$fc.Range("H2").Value = 13.69
$fc.Range("L2").Value = 0.9

# Row 3 (W11)
$fc.Range("D3").Value = 60
$fc.Range("H3").Value = 12.48
$fc.Range("L3").Value = 0.9

# Row 4 (W12)
$fc.Range("D4").Value = 62
$fc.Range("H4").Value = 11.11
$fc.Range("L4").Value = 0.9399999999999999

# Row 5 (W13)
$fc.Range("D5").Value = 64
$fc.Range("H5").Value = 9.800000000000001
$fc.Range("L5").Value = 0.99

# Row 6 (W14)
$fc.Range("D6").Value = 65
$fc.Range("H6").Value = 8.66
$fc.Range("L6").Value = 1.12

# Row 7 (W15)
$fc.Range("D7").Value = 66
$fc.Range("H7").Value = 7.55
$fc.Range("L7").Value = 0.87

# Row 8 (W16)
$fc.Range("D8").Value = 67
$fc.Range("H8").Value = 6.45
$fc.Range("L8").Value = 0.9399999999999999

# Row 9 (W17)
$fc.Range("D9").Value = 68
$fc.Range("H9").Value = 5.37
$fc.Range("L9").Value = 1.07

# Row 10 (W18)
$fc.Range("D10").Value = 70
$fc.Range("H10").Value = 4.24
$fc.Range("L10").Value = 0.89

# Row 11 (W19)
$fc.Range("D11").Value = 71
$fc.Range("H11").Value = 3.2
$fc.Range("L11").Value = 0.96

# Row 12 (W20)
$fc.Range("D12").Value = 72
$fc.Range("H12").Value = 2.17
$fc.Range("L12").Value = 0.92

# Row 13 (W21)
$fc.Range("D13").Value = 73
$fc.Range("H13").Value = 1.15
$fc.Range("L13").Value = 1.03

# Row 14 (W22) - also flips to Urgent/High (4wk low sales check)
$fc.Range("D14").Value = 75
$fc.Range("H14").Value = 0.15
$fc.Range("I14").Value = "High"
$fc.Range("J14").Value = "Urgent"
$fc.Range("L14").Value = 1.12

# Row 15 (W23) - also flips to Urgent/High
$fc.Range("D15").Value = 76
$fc.Range("H15").Value = 0
$fc.Range("I15").Value = "High"
$fc.Range("J15").Value = "Urgent"
$fc.Range("L15").Value = 0.83

# Row 16 (W24) - also flips to Urgent/High
$fc.Range("D16").Value = 78
$fc.Range("H16").Value = 0
$fc.Range("I16").Value = "High"
$fc.Range("J16").Value = "Urgent"
$fc.Range("L16").Value = 1.06

# Row 17 (W25) - also flips to Urgent/High
$fc.Range("D17").Value = 79
$fc.Range("H17").Value = 0
$fc.Range("I17").Value = "High"
$fc.Range("J17").Value = "Urgent"
$fc.Range("L17").Value = 0.9

$sum = $wb.Worksheets.Item("Summary")

# These cells store the figures as text (matching the rest of the Summary
# column), so force a text number format before assigning the string value.
$sum.Range("B9").NumberFormat = "@"
$sum.Range("B9").Value = "1105"

$sum.Range("B10").NumberFormat = "@"
$sum.Range("B10").Value = "511"

$sum.Range("B11").NumberFormat = "@"
$sum.Range("B11").Value = "245"

$sum.Range("B12").NumberFormat = "@"
$sum.Range("B12").Value = "79"

$sum.Range("B14").NumberFormat = "@"
$sum.Range("B14").Value = "59"
